# Actualización automática del tracker
# Adds new match rows (69-75) to the results tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row=69; A=14802841; B="2025-10-01"; C="Laslo Djere";          D="Valentin Vacherot";     E="Gana Laslo Djere";              F=2.5  },
    @{ Row=70; A=14802839; B="2025-10-01"; C="Yoshihito Nishioka";   D="Alexander Shevchenko";  E="Gana Alexander Shevchenko";     F=2.1  },
    @{ Row=71; A=14792374; B="2025-10-01"; C="Zhizhen Zhang";        D="Sebastián Báez";        E="Gana Sebastián Báez";           F=2.3  },
    @{ Row=72; A=14803481; B="2025-10-02"; C="Francisco Comesaña";   D="Ugo Blanchet";          E="Gana Francisco Comesaña";       F=1.73 },
    @{ Row=73; A=14792386; B="2025-10-02"; C="Aleksandar Kovacevic"; D="Juncheng Shang";        E="Gana Aleksandar Kovacevic";     F=3    },
    @{ Row=74; A=14793807; B="2025-10-01"; C="Marta Kostyuk";        D="Jessica Pegula";        E="Gana Jessica Pegula";           F=1.57 },
    @{ Row=75; A=14763285; B="2025-10-01"; C="Mees Rottgering";      D="Mikhail Kukushkin";     E="Gana Mees Rottgering";          F=2.5  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    # Write the match date as literal text (not an Excel date serial):
    # format the cell as Text before assigning, then strip the formatting
    # back off so the stored cell keeps the default (unstyled) look.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    # resultado / profit are still unknown for these new matches -> blank
    # text cells (not "empty"/null), matching the tracker's pending rows.
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).ClearFormats()
    $ws.Cells.Item($row, 8).Value = "'"
    $ws.Cells.Item($row, 8).ClearFormats()
}
